$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-ChoiceHeaderCell($cellRef, $letter) {
    $r = $ws.Range($cellRef)
    $r.Value = "選択肢" + $letter
    $r.Font.Name = "Arial"
    $kanji = $r.Characters(1, 3)
    $kanji.Font.Name = "MS Gothic"
    $romaji = $r.Characters(4, 1)
    $romaji.Font.Name = "Arial"
}

Set-ChoiceHeaderCell "D1" "A"
Set-ChoiceHeaderCell "E1" "B"
Set-ChoiceHeaderCell "F1" "C"

[void]$ws.Range("F2").Select()
